$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New error rows loaded for March data (FF_* feed files).
$data = @(
    @("FF_DailyOrder-26032015.csv",   "incorrect file name",                               ""),
    @("FF_Account-13032015.csv",      "Row 4, Column 52: String data, right truncation",   "Некорректные данные  в поле ClientUserId"),
    @("FF_Account-14032015.csv",      "Row 2, Column 52: String data, right truncation",   "Некорректные данные  в поле ClientUserId"),
    @("FF_Account-18032015.csv",      "Row 15, Column 52: String data, right truncation",  "Некорректные данные  в поле ClientUserId"),
    @("FF_Account-19032015.csv",      "Row 3, Column 52: String data, right truncation",   "Некорректные данные  в поле ClientUserId"),
    @("FF_Order-14032015.csv",        "Row 16, Column 4: String data, right truncation",   "Некорректные данные  в поле ClientUserId"),
    @("FF_Order-19032015.csv",        "Row 3, Column 4: String data, right truncation",    "Некорректные данные  в поле ClientUserId"),
    @("FF_Subscriptions-14032015.csv","Row 1, Column 15: String data, right truncation",   "Некорректные данные  в поле ClientUserId"),
    @("FF_Subscriptions-19032015.csv","Row 1, Column 15: String data, right truncation",   "Некорректные данные  в поле ClientUserId")
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne "") {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
}

# Resize columns to fit the newly loaded content.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("B12").Select() | Out-Null
